$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing B-column values (rows 1-18)
$ws.Range("B1").Value = 10
$ws.Range("B2").Value = 0.01
$ws.Range("B3").Value = 0.02
$ws.Range("B4").Value = 0.5
$ws.Range("B8").Value = 110
$ws.Range("B9").Value = 30
$ws.Range("B11").Value = 720
$ws.Range("B12").Value = 231
$ws.Range("B13").Value = 3
$ws.Range("B14").Value = 103
$ws.Range("B15").Value = 0.05
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 150
$ws.Range("B18").Value = 5

# New row 19: B19 holds an #N/A error value
$ws.Range("B19").Value = "#N/A"

# Update selection to match the new active cell / range
$ws.Range("B1:B19").Select()
